$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Activate() | Out-Null

# Row 13 (SEQ 5 / "長度") column type/length change: VARCHAR2(100) -> NVARCHAR2(200)
$ws.Range("D13").Value = "NVARCHAR2"
$ws.Range("E13").Value = 200

$ws.Range("D13").Select() | Out-Null
